$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 182: add a third cell (C182 = "day 7") to the existing "length" / "长度" row ---
$ws.Range("C182").Value = "day 7"

# --- New rows 191-196: six fresh vocabulary entries, appended after the last used row (190) ---
# Row 191 also carries a "day 8" marker in column C, mirroring the day 7/8 pattern already
# present earlier in the sheet.
$ws.Range("C191").Value = "day 8"
$ws.Range("A191").Value = "portable"
$ws.Range("B191").Value = "便携式的"

$ws.Range("A192").Value = "platform"
$ws.Range("B192").Value = "平台"

$ws.Range("A193").Value = "re-use"
$ws.Range("B193").Value = "再利用"

$ws.Range("A194").Value = "thread"
$ws.Range("B194").Value = "线"

$ws.Range("A195").Value = "assign"
$ws.Range("B195").Value = "分配"

$ws.Range("A196").Value = "maintain"
$ws.Range("B196").Value = "保持"

# --- Update the sheet's saved selection/scroll state to match where the author ended up editing ---
$ws.Range("B200").Select()
